$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 1066, pushing existing data (old rows 1066-1113)
# down to become rows 1069-1116.
$ws.Range("A1066:A1068").EntireRow.Insert()

# Common template values shared by every row in this block (unchanged by the edit).
$company = 'Agrícola del Norte S.A. de Arica'
$region = 'Arica y Parinacota'
$codreg = 15
$catId = 100112020
$categoria = 'Tomate'
$variedad = 'Larga vida'
$unidad = '$/caja 10 kilos'
$origen = 'Región de Arica y Parinacota'
$kg = 10
$clasificacion = 'Hortaliza'

# New data for the three inserted rows (row, calidad, fecha, volumen, pmin, pmax, pprom, pkg)
$newRows = @(
  @(1066, 'Primera', 44826, 250, 2000, 2500, 2250, 225),
  @(1067, 'Segunda', 44826, 300, 1500, 2000, 1750, 175),
  @(1068, 'Tercera', 44826, 300, 1000, 1500, 1250, 125)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $calidad = $r[1]
    $fecha = $r[2]
    $volumen = $r[3]
    $pmin = $r[4]
    $pmax = $r[5]
    $pprom = $r[6]
    $pkg = $r[7]

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $company
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $catId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $pmin
    $ws.Cells.Item($row, 12).Value = $pmax
    $ws.Cells.Item($row, 13).Value = $pprom
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $pkg
    $ws.Cells.Item($row, 17).Value = $kg
    $ws.Cells.Item($row, 18).Value = $clasificacion

    # Ensure the date cell uses the same date style as the rest of column D.
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 3, 4).NumberFormat
}

Write-Host "Done inserting new rows."
